$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D contains price strings that must stay as text (not auto-converted
# to numbers by Excel), matching the original inlineStr cell formatting.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.104.79"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.655.21"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.76"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5271"
$ws.Range("E6").Value = "  +2.16%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2611"
$ws.Range("E8").Value = "  -0.91%  "
$ws.Range("E9").Value = "  +1.22%  "
$ws.Range("E10").Value = "  -1.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07793"
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.508"
$ws.Range("E12").Value = "  +1.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.597.09"
$ws.Range("E13").Value = "  -3.68%  "
$ws.Range("E14").Value = "  +1.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅8230"
$ws.Range("E15").Value = "  +1.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.42"
$ws.Range("E16").Value = "  +0.98%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.117.39"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.584"
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "190.87"
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.034"
$ws.Range("E22").Value = "  +0.60%  "
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "141.51"
$ws.Range("E24").Value = "  +1.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1237"
$ws.Range("E25").Value = "  +1.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.245"
$ws.Range("E26").Value = "  +0.83%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.429"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05892"
$ws.Range("E29").Value = "  -1.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.275"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.520"
$ws.Range("E31").Value = "  -1.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.264"
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("E33").Value = "  -0.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9525"
$ws.Range("E34").Value = "  -1.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.784"
$ws.Range("E35").Value = "  +0.45%  "
$ws.Range("E36").Value = "  -0.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5702"
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("E38").Value = "  +1.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.803"
$ws.Range("E39").Value = "  -2.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8483"
$ws.Range("E40").Value = "  -0.73%  "
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.028.67"
$ws.Range("E42").Value = "  +2.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "102.59"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.799.60"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "57.12"
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  -0.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4298"
$ws.Range("E47").Value = "  +2.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.477"
$ws.Range("E48").Value = "  +2.08%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.853"
$ws.Range("E49").Value = "  -1.93%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05154"
$ws.Range("E50").Value = "  -0.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.09692"
$ws.Range("E51").Value = "  -0.44%  "
